$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Character Data sheet - append new character rows (3-10)
# ---------------------------------------------------------------------------
$charData = $wb.Worksheets.Item("Character Data")

$charRows = @(
    @("Super Smash Bros Melee", "Mario", "Super Mario Bros", 18, 5),
    @("Super Smash Bros Melee", "Bowser", "Super Mario Bros", 19, 5),
    @("Super Smash Bros Melee", "Peach", "Super Mario Bros", 15, 5),
    @("Super Smash Bros Ultimate", "Yoshi", "Yoshi", 20, 5),
    @("Super Smash Bros Ultimate", "Mario", "Super Mario Bros", 18, 5),
    @("Super Smash Bros Ultimate", "Bowser", "Super Mario Bros", 19, 5),
    @("Super Smash Bros Ultimate", "Peach", "Super Mario Bros", 15, 5),
    @("Super Smash Bros Ultimate", "Duck Hunt", "Duck Hunt", 15, 5)
)

$r = 3
foreach ($row in $charRows) {
    $charData.Cells.Item($r, 1).Value = $row[0]
    $charData.Cells.Item($r, 2).Value = $row[1]
    $charData.Cells.Item($r, 3).Value = $row[2]
    $charData.Cells.Item($r, 4).Value = $row[3]
    $charData.Cells.Item($r, 5).Value = $row[4]
    $r++
}

$charData.Range("D12").Select()

# ---------------------------------------------------------------------------
# 2. Move Sheet - fix character names on existing rows + append new moves
# ---------------------------------------------------------------------------
$moveSheet = $wb.Worksheets.Item("Move Sheet")

$moveSheet.Range("A3").Value = "Bowser"
$moveSheet.Range("A4").Value = "Mario"
$moveSheet.Range("A5").Value = "Mario"
$moveSheet.Range("C5").Value = "strong"

# E5/F5 are formatted as Text ("@"); round-trip through the Normal style so
# the numeric literal is stored as a real number (matching the target) while
# putting the original Text number format back on the cell afterwards.
$moveSheet.Range("E5").Style = "Normal"
$moveSheet.Range("E5").Value = 9
$moveSheet.Range("E5").NumberFormat = "@"
$moveSheet.Range("F5").Style = "Normal"
$moveSheet.Range("F5").Value = 11
$moveSheet.Range("F5").NumberFormat = "@"

$newMoveRows = @(
    @("Yoshi", "Super Smash Bros Ultimate", "smash", "up", 10, 15),
    @("Bowser", "Super Smash Bros Ultimate", "special", "down", 12, 13),
    @("Mario", "Super Smash Bros Ultimate", "strong", "up", 9, 11),
    @("Mario", "Super Smash Bros Ultimate", "strong", "down", 9, 11)
)

$r = 6
foreach ($row in $newMoveRows) {
    $moveSheet.Cells.Item($r, 1).Value = $row[0]
    $moveSheet.Cells.Item($r, 2).Value = $row[1]
    $moveSheet.Cells.Item($r, 3).Value = $row[2]
    $moveSheet.Cells.Item($r, 4).Value = $row[3]

    # DamageStartFrame/DamageEndFrame columns carry on the Text ("@") number
    # format used by the rows above, so round-trip through Normal again.
    $eCell = $moveSheet.Cells.Item($r, 5)
    $eCell.Style = "Normal"
    $eCell.Value = $row[4]
    $eCell.NumberFormat = "@"

    $fCell = $moveSheet.Cells.Item($r, 6)
    $fCell.Style = "Normal"
    $fCell.Value = $row[5]
    $fCell.NumberFormat = "@"

    $r++
}

$moveSheet.Range("E12").Select()

# ---------------------------------------------------------------------------
# 3. Stage Data - append new stage row (12)
# ---------------------------------------------------------------------------
$stageData = $wb.Worksheets.Item("Stage Data")

$stageData.Range("A12").Value = "Duck Hunt"
$stageData.Range("B12").Value = "Duck Hunt"
$stageData.Range("C12").Value = "Duck Hunt"
$stageData.Range("D12").Value = "Super Smash Bros Ultimate"

$stageData.Range("A12:D12").Select()

# ---------------------------------------------------------------------------
# 4. Item Data - append new item row (7)
# ---------------------------------------------------------------------------
$itemData = $wb.Worksheets.Item("Item Data")

$itemData.Range("A7").Value = "Hammer"
$itemData.Range("B7").Value = "Super Smash Bros."
$itemData.Range("C7").Value = 1
$itemData.Range("D7").Value = "Super Smash Bros Ultimate"

$itemData.Range("B10").Select()

# ---------------------------------------------------------------------------
# 5. New "User" sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$userSheet = $wb.Worksheets.Add($null, $lastSheet)
$userSheet.Name = "User"

$userSheet.Columns.Item(1).ColumnWidth = 16.6640625
$userSheet.Columns.Item(2).ColumnWidth = 16.77734375

$userSheet.Range("A1").Value = "UserName"
$userSheet.Range("B1").Value = "Password"
$userSheet.Range("A1:B1").Style = "Neutral"

$userRows = @(
    @("mayor", "mayor"),
    @("test", "team"),
    @("answer", "here")
)

$r = 2
foreach ($row in $userRows) {
    $userSheet.Cells.Item($r, 1).Value = $row[0]
    $userSheet.Cells.Item($r, 2).Value = $row[1]
    $r++
}

$userSheet.Range("A1:B1").Select()

# ---------------------------------------------------------------------------
# 6. New "UserFavorites" sheet
# ---------------------------------------------------------------------------
$favSheet = $wb.Worksheets.Add($null, $userSheet)
$favSheet.Name = "UserFavorites"

$favSheet.Range("A1").Value = "UserName"
$favSheet.Range("B1").Value = "Game"
$favSheet.Range("C1").Value = "Character"
$favSheet.Range("A1:C1").Style = "Neutral"

$favSheet.Range("A2").Value = "mayor"
$favSheet.Range("B2").Value = "Super Smash Bros Melee"
$favSheet.Range("C2").Value = "Mario"

$favSheet.Range("A3").Value = "test"
$favSheet.Range("B3").Value = "Super Smash Bros Ultimate"
$favSheet.Range("C3").Value = "Bowser"

$favSheet.Range("D11").Select()
$favSheet.Activate()
